$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("I9").Value = "1402-04-14 (9)"
$ws.Range("M9").Value = "1402-04-14 (2)"

$ws.Range("I12").Value = -8031256
$ws.Range("M12").Value = 28723842

$ws.Range("I14").Value = -8046676
$ws.Range("M14").Value = 28723842

$ws.Range("I16").Value = 27181
$ws.Range("M16").Value = 18730

$ws.Range("M17").Value = -18837232

$ws.Range("I29").Value = -24836
$ws.Range("M29").Value = -39170

$ws.Range("I32").Value = -1461406
$ws.Range("M32").Value = -19223782

$ws.Range("E36").Value = 0
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 0
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 0
